$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "rpZLs656"
$ws.Range("B2").Value = 23080424
$ws.Range("C2").Value = "hkmkdsh13"
$ws.Range("D2").Value = "Dma5W7%!"
$ws.Range("F2").Value = "NtQlFPxc"
$ws.Range("G2").Value = "WQWc"
